$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CLAVES")

$ws.Range("A4").Value = "maria@reqres.in"
$ws.Range("B4").Value = "maria"

$ws.Range("A6").Select()
